$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2:E2').NumberFormat = '@'
$ws.Range('D2').Value = '24.883.46'
$ws.Range('E2').Value = '  +1.93%  '

$ws.Range('D3:E3').NumberFormat = '@'
$ws.Range('D3').Value = '1.710.17'
$ws.Range('E3').Value = '  +1.83%  '

$ws.Range('D4:E4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5:E5').NumberFormat = '@'
$ws.Range('D5').Value = '311.16'
$ws.Range('E5').Value = '  +1.45%  '

$ws.Range('D6:E6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  +0.30%  '

$ws.Range('D7:E7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3753'
$ws.Range('E7').Value = '  +1.23%  '

$ws.Range('D8:E8').NumberFormat = '@'
$ws.Range('D8').Value = '49.67'
$ws.Range('E8').Value = '  +3.07%  '

$ws.Range('D9:E9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3442'
$ws.Range('E9').Value = '  +0.02%  '

$ws.Range('D10:E10').NumberFormat = '@'
$ws.Range('D10').Value = '1.207'
$ws.Range('E10').Value = '  +2.09%  '

$ws.Range('D11:E11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07535'
$ws.Range('E11').Value = '  +3.65%  '

$ws.Range('D12:E12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9998'
$ws.Range('E12').Value = '  +0.06%  '

$ws.Range('D13:E13').NumberFormat = '@'
$ws.Range('D13').Value = '21.09'
$ws.Range('E13').Value = '  +3.20%  '

$ws.Range('D14:E14').NumberFormat = '@'
$ws.Range('D14').Value = '6.294'
$ws.Range('E14').Value = '  +2.99%  '

$ws.Range('D15:E15').NumberFormat = '@'
$ws.Range('D15').Value = '7.036'
$ws.Range('E15').Value = '  +4.26%  '

$ws.Range('D16:E16').NumberFormat = '@'
$ws.Range('D16').Value = '1.706.91'
$ws.Range('E16').Value = '  +1.78%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001136'

$ws.Range('D18:E18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06732'
$ws.Range('E18').Value = '  +0.16%  '

$ws.Range('D19:E19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  +0.30%  '

$ws.Range('D20:E20').NumberFormat = '@'
$ws.Range('D20').Value = '84.54'
$ws.Range('E20').Value = '  +4.18%  '

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +5.36%  '

$ws.Range('D22:E22').NumberFormat = '@'
$ws.Range('D22').Value = '6.378'
$ws.Range('E22').Value = '  +4.47%  '

$ws.Range('D23:E23').NumberFormat = '@'
$ws.Range('D23').Value = '13.22'
$ws.Range('E23').Value = '  +10.54%  '

$ws.Range('D24:E24').NumberFormat = '@'
$ws.Range('D24').Value = '24.821.74'
$ws.Range('E24').Value = '  +1.92%  '

$ws.Range('D25:E25').NumberFormat = '@'
$ws.Range('D25').Value = '2.445'
$ws.Range('E25').Value = '  +0.49%  '

$ws.Range('D26:E26').NumberFormat = '@'
$ws.Range('D26').Value = '2.786'
$ws.Range('E26').Value = '  +4.42%  '

$ws.Range('D27:E27').NumberFormat = '@'
$ws.Range('D27').Value = '20.38'
$ws.Range('E27').Value = '  +3.99%  '

$ws.Range('D28:E28').NumberFormat = '@'
$ws.Range('D28').Value = '151.95'
$ws.Range('E28').Value = '  -0.28%  '

$ws.Range('D29:E29').NumberFormat = '@'
$ws.Range('D29').Value = '132.27'
$ws.Range('E29').Value = '  +3.88%  '

$ws.Range('D30:E30').NumberFormat = '@'
$ws.Range('D30').Value = '1.901.07'
$ws.Range('E30').Value = '  +2.14%  '

$ws.Range('D31:E31').NumberFormat = '@'
$ws.Range('D31').Value = '1.237'
$ws.Range('E31').Value = '  +27.39%  '

$ws.Range('D32:E32').NumberFormat = '@'
$ws.Range('D32').Value = '6.903'
$ws.Range('E32').Value = '  +9.22%  '

$ws.Range('D33:E33').NumberFormat = '@'
$ws.Range('D33').Value = '4.253'
$ws.Range('E33').Value = '  +5.69%  '

$ws.Range('D34:E34').NumberFormat = '@'
$ws.Range('D34').Value = '1.823'
$ws.Range('E34').Value = '  +5.01%  '

$ws.Range('D35:E35').NumberFormat = '@'
$ws.Range('D35').Value = '13.75'
$ws.Range('E35').Value = '  +11.53%  '

$ws.Range('D36:E36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08773'
$ws.Range('E36').Value = '  +3.42%  '

$ws.Range('D37:E37').NumberFormat = '@'
$ws.Range('D37').Value = '5.610'
$ws.Range('E37').Value = '  +4.93%  '

$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38:E38').NumberFormat = '@'
$ws.Range('D38').Value = '9.319'
$ws.Range('E38').Value = '  +1.98%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39:E39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06695'
$ws.Range('E39').Value = '  +2.95%  '

$ws.Range('D40:E40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02414'
$ws.Range('E40').Value = '  +3.25%  '

$ws.Range('D41:E41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2237'
$ws.Range('E41').Value = '  +5.84%  '

$ws.Range('D42:E42').NumberFormat = '@'
$ws.Range('D42').Value = '1.277'
$ws.Range('E42').Value = '  +1.26%  '

$ws.Range('D43:E43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6445'
$ws.Range('E43').Value = '  +4.16%  '

$ws.Range('D44:E44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9993'
$ws.Range('E44').Value = '  +0.27%  '

$ws.Range('D45:E45').NumberFormat = '@'
$ws.Range('D45').Value = '13.92'
$ws.Range('E45').Value = '  +7.34%  '

$ws.Range('D46:E46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6168'
$ws.Range('E46').Value = '  +3.60%  '

$ws.Range('D47:E47').NumberFormat = '@'
$ws.Range('D47').Value = '3.829'
$ws.Range('E47').Value = '  +1.27%  '

$ws.Range('D48:E48').NumberFormat = '@'
$ws.Range('D48').Value = '2.139'
$ws.Range('E48').Value = '  +5.55%  '

$ws.Range('D49:E49').NumberFormat = '@'
$ws.Range('D49').Value = '130.18'
$ws.Range('E49').Value = '  +2.35%  '

$ws.Range('D50:E50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07313'
$ws.Range('E50').Value = '  +1.28%  '

$ws.Range('D51:E51').NumberFormat = '@'
$ws.Range('D51').Value = '79.74'
$ws.Range('E51').Value = '  +5.16%  '
